$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12: date label as plain text (avoid Excel's automatic date parsing)
# so it matches the existing rows' shared-string-text representation.
$ws.Range("A12").Value = "'2025-08-30"
$ws.Range("A12").Style = "Normal"

# Refreshed nowcast figures for existing rows 2-11, plus the new row 12 values.
$ws.Range("B2").Value = [double]"0.35732963346269386"
$ws.Range("C2").Value = [double]"0"
$ws.Range("D2").Value = [double]"0"
$ws.Range("E2").Value = [double]"0"
$ws.Range("F2").Value = [double]"0"
$ws.Range("G2").Value = [double]"0"
$ws.Range("H2").Value = [double]"0"
$ws.Range("I2").Value = [double]"0"
$ws.Range("J2").Value = [double]"0"
$ws.Range("K2").Value = [double]"0"
$ws.Range("B3").Value = [double]"0.32518399396045433"
$ws.Range("C3").Value = [double]"0"
$ws.Range("D3").Value = [double]"-0.005553618388527283"
$ws.Range("E3").Value = [double]"6.630278573450277e-05"
$ws.Range("F3").Value = [double]"-0.0013606276388917259"
$ws.Range("G3").Value = [double]"0.00010394158051328925"
$ws.Range("H3").Value = [double]"-9.886914845301429e-05"
$ws.Range("I3").Value = [double]"-0.001121538971848243"
$ws.Range("J3").Value = [double]"0"
$ws.Range("K3").Value = [double]"0.005509764094738079"
$ws.Range("B4").Value = [double]"0.3246823271591469"
$ws.Range("C4").Value = [double]"-0.0025823633850423472"
$ws.Range("D4").Value = [double]"0"
$ws.Range("E4").Value = [double]"0.0005402635467678584"
$ws.Range("F4").Value = [double]"7.591562601938678e-05"
$ws.Range("G4").Value = [double]"0"
$ws.Range("H4").Value = [double]"0.00017815094292701625"
$ws.Range("I4").Value = [double]"-0.0024784699503127947"
$ws.Range("J4").Value = [double]"0.0005359048495616577"
$ws.Range("K4").Value = [double]"0.00032248442744131633"
$ws.Range("B5").Value = [double]"0.3646698119966691"
$ws.Range("C5").Value = [double]"0.010582493435496831"
$ws.Range("D5").Value = [double]"-0.007580929477483236"
$ws.Range("E5").Value = [double]"0.00030217341433305347"
$ws.Range("F5").Value = [double]"0.0010451204854101351"
$ws.Range("G5").Value = [double]"-0.0008829490443779485"
$ws.Range("H5").Value = [double]"3.6631627045577925e-05"
$ws.Range("I5").Value = [double]"-0.0008553640916597865"
$ws.Range("J5").Value = [double]"0"
$ws.Range("K5").Value = [double]"-0.0031210652726189614"
$ws.Range("B6").Value = [double]"0.39171244027192104"
$ws.Range("C6").Value = [double]"0.027504414683785415"
$ws.Range("D6").Value = [double]"0"
$ws.Range("E6").Value = [double]"-0.000292700031868534"
$ws.Range("F6").Value = [double]"2.344017861949787e-05"
$ws.Range("G6").Value = [double]"0"
$ws.Range("H6").Value = [double]"-8.077501065895375e-05"
$ws.Range("I6").Value = [double]"-0.0019070943150054553"
$ws.Range("J6").Value = [double]"0"
$ws.Range("K6").Value = [double]"-0.005946532576799002"
$ws.Range("B7").Value = [double]"0.31892492538086936"
$ws.Range("C7").Value = [double]"0"
$ws.Range("D7").Value = [double]"-0.0031521999452441365"
$ws.Range("E7").Value = [double]"-0.0015451521624765732"
$ws.Range("F7").Value = [double]"-0.007613838283928674"
$ws.Range("G7").Value = [double]"0.001023098162061924"
$ws.Range("H7").Value = [double]"0"
$ws.Range("I7").Value = [double]"0.00026292891325876657"
$ws.Range("J7").Value = [double]"0"
$ws.Range("K7").Value = [double]"-0.0032753170484974836"
$ws.Range("B8").Value = [double]"0.2250492513677827"
$ws.Range("C8").Value = [double]"-0.060527042791817594"
$ws.Range("D8").Value = [double]"0"
$ws.Range("E8").Value = [double]"-0.0001470932620413092"
$ws.Range("F8").Value = [double]"-0.00032393204950123734"
$ws.Range("G8").Value = [double]"0"
$ws.Range("H8").Value = [double]"6.226863214096356e-05"
$ws.Range("I8").Value = [double]"0.0015898728951537263"
$ws.Range("J8").Value = [double]"0"
$ws.Range("K8").Value = [double]"-0.0007237550799384929"
$ws.Range("B9").Value = [double]"0.28821937053307406"
$ws.Range("C9").Value = [double]"0"
$ws.Range("D9").Value = [double]"0.0061824108848021"
$ws.Range("E9").Value = [double]"-0.0031074043217820656"
$ws.Range("F9").Value = [double]"-0.004634646756981129"
$ws.Range("G9").Value = [double]"0.002028293762944151"
$ws.Range("H9").Value = [double]"-0.0001898400656179828"
$ws.Range("I9").Value = [double]"-7.177135270133911e-05"
$ws.Range("J9").Value = [double]"0"
$ws.Range("K9").Value = [double]"0.0025738345080049407"
$ws.Range("B10").Value = [double]"0.4378941889091869"
$ws.Range("C10").Value = [double]"0.10202457512007569"
$ws.Range("D10").Value = [double]"0"
$ws.Range("E10").Value = [double]"-0.0009057822040697325"
$ws.Range("F10").Value = [double]"-0.0001929331229023158"
$ws.Range("G10").Value = [double]"0"
$ws.Range("H10").Value = [double]"-5.631698127426953e-06"
$ws.Range("I10").Value = [double]"8.06278417066223e-05"
$ws.Range("J10").Value = [double]"-0.0023025253424738285"
$ws.Range("K10").Value = [double]"0.00441474893636945"
$ws.Range("B11").Value = [double]"0.2953293387166088"
$ws.Range("C11").Value = [double]"0"
$ws.Range("D11").Value = [double]"-0.034585966520740495"
$ws.Range("E11").Value = [double]"0.0032252475482832626"
$ws.Range("F11").Value = [double]"0.006089927290486877"
$ws.Range("G11").Value = [double]"0.0031451476749442"
$ws.Range("H11").Value = [double]"0.0010889405499454228"
$ws.Range("I11").Value = [double]"0.0016644051435757475"
$ws.Range("J11").Value = [double]"0"
$ws.Range("K11").Value = [double]"-0.004092206017015043"
$ws.Range("B12").Value = [double]"0.2370759115919914"
$ws.Range("C12").Value = [double]"-0.06924194231654604"
$ws.Range("D12").Value = [double]"0"
$ws.Range("E12").Value = [double]"0.002553578605936005"
$ws.Range("F12").Value = [double]"4.7492376585786696e-05"
$ws.Range("G12").Value = [double]"0"
$ws.Range("H12").Value = [double]"9.411632500719862e-06"
$ws.Range("I12").Value = [double]"-0.0017449468792976681"
$ws.Range("J12").Value = [double]"0"
$ws.Range("K12").Value = [double]"-0.018327409138999995"
